# Ejemplo2.xlsx - "cocina output process"
# Updates the reference/quantity table on Hoja1 and the single id/value
# pair on Hoja2 with the new "cocina" data set, formats the recalculated
# quantity in B11 with a thousands-separator number format, and leaves the
# selection where the author last clicked before saving.

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item(1)   # Hoja1
$ws2 = $wb.Worksheets.Item(2)   # Hoja2

# --- Hoja1: Ref / Cantidad table -------------------------------------------
# Row 1 (headers "Ref" / "Cantidad") is unchanged.

$ws1.Range("A2").Value = "RH150010"
$ws1.Range("B2").Value = 193

$ws1.Range("A3").Value = "G4504"
$ws1.Range("B3").Value = 1

$ws1.Range("A4").Value = "E2671"
$ws1.Range("B4").Value = 1

$ws1.Range("A5").Value = "AIN0L10"
$ws1.Range("B5").Value = 1

$ws1.Range("A6").Value = "M0F015"
$ws1.Range("B6").Value = 28

$ws1.Range("A7").Value = "T1502015"
$ws1.Range("B7").Value = 72

$ws1.Range("A8").Value = "TS03981"
$ws1.Range("B8").Value = 60

$ws1.Range("A9").Value = "BSCO0051"
$ws1.Range("B9").Value = 30

$ws1.Range("A10").Value = "SP2"
$ws1.Range("B10").Value = 88

$ws1.Range("A11").Value = "MG023"
$ws1.Range("B11").NumberFormat = "#,##0"
$ws1.Range("B11").Value = 69

$ws1.Range("A12").Value = "#03"
$ws1.Range("B12").Value = 1

# Match the default column width the author's workbook ended up with.
$ws1.StandardWidth = 11.5546875

$ws1.Activate()
$ws1.Range("E8").Select()

# --- Hoja2: single id/value pair -------------------------------------------

$ws2.Range("A1").Value = "id"
$ws2.Range("A2").Value = 1

$ws2.StandardWidth = 11.5546875

$ws2.Activate()
$ws2.Range("D8").Select()

# Leave Hoja1 as the active/selected sheet, as in the saved workbook.
$ws1.Activate()
